$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The "_GoBack" bookmark currently sits as an (empty) bookmarkStart/
#    bookmarkEnd pair right after the closing "}" of the last code
#    sample in the document. Word always re-homes "_GoBack" to the
#    most recent edit location, so the stale one is removed here --
#    it gets re-created below, around the paragraph that is actually
#    edited.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the "Eloltesztelos ciklus" heading paragraph and make its
#    text bold.
# ------------------------------------------------------------------
$headingPara = $null

$findRange = $d.Content
$found = $findRange.Find.Execute("Elöltesztelős ciklus", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $headingPara = $findRange.Paragraphs(1)
}

if ($headingPara -eq $null) {
    # Fallback: scan paragraphs directly for the exact heading text.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs($i)
        if ($candidate.Range.Text -eq "Elöltesztelős ciklus`r") {
            $headingPara = $candidate
            break
        }
    }
}

$headingPara.Range.Font.Bold = 1

# ------------------------------------------------------------------
# 3) Re-create "_GoBack" spanning the whole paragraph that was just
#    edited (mirroring Word's own behaviour after an in-place edit).
# ------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $headingPara.Range)
